$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Point Anomalies (Validation noise added) - Univariate HMM (row 5)
$ws.Range("Q5").Value = 59064
$ws.Range("R5").Value = 31048
$ws.Range("S5").Value = 23239
$ws.Range("T5").Value = 17031
$ws.Range("U5").Value = 11585
$ws.Range("V5").Value = 9246

# Point Anomalies (Validation noise added) - Multivariate HMM (row 8)
$ws.Range("Q8").Value = 57438
$ws.Range("R8").Value = 34935
$ws.Range("S8").Value = 24086
$ws.Range("T8").Value = 16903
$ws.Range("U8").Value = 14571
$ws.Range("V8").Value = 12960

# Point Anomalies (Validation noise added) - Multivariate HMM (row 10)
$ws.Range("Q10").Value = 54919
$ws.Range("R10").Value = 30045
$ws.Range("S10").Value = 23262
$ws.Range("T10").Value = 19491
$ws.Range("U10").Value = 16167
$ws.Range("V10").Value = 12532

# Collective Anomalies (Validation noise added) - Univariate HMM (rows 17-19)
$ws.Range("Q17").Value = 14339
$ws.Range("R17").Value = 6503
$ws.Range("S17").Value = 4743
$ws.Range("T17").Value = 3433

$ws.Range("Q18").Value = 10307
$ws.Range("R18").Value = 5149
$ws.Range("S18").Value = 4061
$ws.Range("T18").Value = 2958

$ws.Range("Q19").Value = 10520
$ws.Range("R19").Value = 5240
$ws.Range("S19").Value = 4028
$ws.Range("T19").Value = 3006

# Collective Anomalies (Validation noise added) - Multivariate HMM (rows 20-22)
$ws.Range("Q20").Value = 14266
$ws.Range("R20").Value = 7101
$ws.Range("S20").Value = 4913
$ws.Range("T20").Value = 3578

$ws.Range("Q21").Value = 12835
$ws.Range("R21").Value = 6174
$ws.Range("S21").Value = 4046
$ws.Range("T21").Value = 3331

$ws.Range("Q22").Value = 13046
$ws.Range("R22").Value = 6267
$ws.Range("S22").Value = 4595
$ws.Range("T22").Value = 3374

# Move the active selection to match the author's final cursor position
$ws.Range("Q9").Select()
